$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carrier prefix values to write into column B (Known Prefixes), keyed by row.
$values = @{
    4  = "ACLU"
    6  = "AMCU, APDU,ANNU, APHU, APIU, APLU, APRU, APZU, CGHU, CGMU, CGTU, CMAU, CMNU, CNCU, DVRU, ECMU, MMCU, NEPU, NOLU, NOSU, NUSU, OPDU, OTAU, STMU"
    7  = "CBHU, CCLU, CSLU, CSNU"
    8  = "CMCU, SEFU"
    11 = "EGHU, EHSU, EISU, EMCU, HMCU, IMTU, LTIU, UGMU"
    12 = "CADU, CNIU, ENAU, GRIU, HASU, KHJU, KHLU, SUDU"
    14 = "AZLU, CASU, CMUU, CPSU, CSQU, CSVU, FANU, HAMU, HLBU, HLCU, HLXU, ITAU, IVLU, LBIU, LNXU, LYKU, MOMU, QIBU, QNNU, TLEU, TMMU, UACU, UAEU, UASU"
    18 = "ICCU"
    19 = "BMLU, KOSU"
    20 = "APMU, CADU, CNIU, COZU, FAAU, FRLU, KNLU, LOTU, LOTU, MAEU, MALU, MCAU, MCHU, MCRU, MHHU, MIEU, MMAU, MNBU, MRKU, MRSU, MSAU, MSFU, MSKU, MSWU, MVIU, MWCU, MWMU, OCLU, POCU, PONU, SCMU, TORU"
    21 = "CXCU, HRZU, MATU"
    23 = "GTIU, MEDU, MSCU, MSDU, MSMU, MSPU, MSYU, MSZU"
    25 = "NSAU"
    26 = "AKLU, EKLU, ESSU, KKFU, KKLU, KKTU, KLFU, KLTU, KXTU, MOAU, MOEU, MOFU, MOGU, MOLU, MORU, MOSU, MOTU, NYKU, ONEU, PXCU"
    27 = "OOCU, OOLU"
    28 = "PCIU, PILU"
    31 = "SEAU"
    33 = "SMLU"
    35 = "SMCU"
    37 = "AWSU"
    38 = "NPRU, SBGU, STRU"
    40 = "SACU, WLNU, WWLU"
    41 = "TPCU, WHLU, WHSU"
    42 = "YMLU, YMMU"
    43 = "ZCLU, ZCSU, ZIMU, ZMOU"
}

# Rows whose B cell did not previously exist in the sheet and therefore need
# the column's default "Text" number format applied explicitly so the cell
# carries style index 1 (numFmtId 49 / "@") when written out.
$needsTextFormat = @(6, 7, 8, 35, 37, 38, 40, 41, 42, 43)

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $values[$row]
    if ($needsTextFormat -contains $row) {
        $cell.NumberFormat = "@"
    }
}

# Restore the previous scroll/selection state as closely as possible and move
# the active selection to match the saved view.
$ws.Range("B14").Select()
